$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2 becomes a numeric value instead of text
$ws.Range("B2").Value = 9876543210

# New row 3 - phone numbers keep their leading zero, so force them as text
$ws.Range("A3").Value = "neil"
$ws.Range("B3").Value = "'0987654321"
$ws.Range("C3").Value = "hyderabad"

# New row 4 (duplicate of row 3)
$ws.Range("A4").Value = "neil"
$ws.Range("B4").Value = "'0987654321"
$ws.Range("C4").Value = "hyderabad"

# New row 5
$ws.Range("A5").Value = "tedst"
$ws.Range("B5").Value = 1234567890
$ws.Range("C5").Value = "delhi"
